$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.941.38'
$ws.Range("E2").Value = '  -1.32%  '

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.430.05'
$ws.Range("E3").Value = '  -0.92%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.40%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.91'
$ws.Range("E5").Value = '  -3.25%  '

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.41'
$ws.Range("E6").Value = '  +13.11%  '

# Row 7 - was XRP, now LidoStakedEther
$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.424.57'
$ws.Range("E7").Value = '  -0.20%  '

# Row 8 - was LidoStakedEther, now XRP
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +1.44%  '

# Row 9 - USDC
$ws.Range("E9").Value = '  -0.35%  '

# Row 10 - Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.678'
$ws.Range("E10").Value = '  +4.19%  '

# Row 11 - Dogecoin
$ws.Range("E11").Value = '  +12.88%  '

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.40'
$ws.Range("E12").Value = '  +5.04%  '

# Row 13 - TRON
$ws.Range("E13").Value = '  -1.40%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.983.02'
$ws.Range("E14").Value = '  -0.95%  '

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.60'
$ws.Range("E15").Value = '  +1.17%  '

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.86'
$ws.Range("E16").Value = '  -0.46%  '

# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.441.96'
$ws.Range("E17").Value = '  -2.17%  '

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.018.58'
$ws.Range("E18").Value = '  -1.59%  '

# Row 19 - Uniswap
$ws.Range("E19").Value = '  +5.36%  '

# Row 20 - Polygon
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.02'
$ws.Range("E20").Value = '  -2.47%  '

# Row 21 - ShibaInu
$ws.Range("E21").Value = '  +15.18%  '

# Row 22 - ImmutableX
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.25'
$ws.Range("E22").Value = '  -4.15%  '

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.50'
$ws.Range("E23").Value = '  +9.26%  '

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.01'
$ws.Range("E24").Value = '  -0.90%  '

# Row 25 - BitcoinCash
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '308.95'
$ws.Range("E25").Value = '  +1.19%  '

# Row 26 - PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.17'
$ws.Range("E26").Value = '  -4.13%  '

# Row 27 - Filecoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.58'
$ws.Range("E27").Value = '  +8.27%  '

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '29.77'
$ws.Range("E28").Value = '  -1.72%  '

# Row 29 - LEO
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.53'
$ws.Range("E29").Value = '  +0.27%  '

# Row 30 - RenderToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.50'
$ws.Range("E30").Value = '  -0.29%  '

# Row 31 - Kaspa
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.173'
$ws.Range("E31").Value = '  -2.21%  '

# Row 32 - Hedera
$ws.Range("E32").Value = '  +1.31%  '

# Row 33 - was InjectiveProtocol, now Cosmos
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.68'
$ws.Range("E33").Value = '  +1.41%  '

# Row 34 - was Cosmos, now InjectiveProtocol
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.98'
$ws.Range("E34").Value = '  +8.09%  '

# Row 35 - Toncoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.56'
$ws.Range("E35").Value = '  +2.57%  '

# Row 37 - VeChain
$ws.Range("E37").Value = '  -4.44%  '

# Row 38 - OKB
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.32'
$ws.Range("E38").Value = '  +0.44%  '

# Row 39 - FirstDigitalUSD
$ws.Range("E39").Value = '  -0.10%  '

# Row 40 - LidoDAOToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.41'
$ws.Range("E40").Value = '  -0.55%  '

# Row 41 - Stacks
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.01'
$ws.Range("E41").Value = '  -4.05%  '

# Row 42 - Stellar
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.125'
$ws.Range("E42").Value = '  +1.49%  '

# Row 43 - ARBITRUM
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.97'
$ws.Range("E43").Value = '  +1.27%  '

# Row 44 - Monero
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '136.29'
$ws.Range("E44").Value = '  -1.87%  '

# Row 45 - TheGraph
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.286'
$ws.Range("E45").Value = '  -1.17%  '

# Row 46 - Celestia
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.96'
$ws.Range("E46").Value = '  +0.64%  '

# Row 47 - NEARProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.91'
$ws.Range("E47").Value = '  -1.78%  '

# Row 48 - WEMIXToken
$ws.Range("E48").Value = '  -2.08%  '

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.68'
$ws.Range("E49").Value = '  -4.13%  '

# Row 50 - RocketPoolETH
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.773.53'
$ws.Range("E50").Value = '  -0.55%  '

# Row 51 - Maker
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.159.46'
$ws.Range("E51").Value = '  -1.94%  '
